{"js": "// Office.js (Word JavaScript API) edit script.\n// Body: async (context) => { ... }\n//\n// The document is a \"two-digit division\" worksheet: a centered date title\n// paragraph followed by a 20-row x 5-column table where every 4th row (rows\n// 0, 4, 8, 12, 16 - zero based) holds a division problem/answer string like\n// \"85\u00f75=17, 0\" and the other rows are blank spacer rows. This script updates\n// the date title and every populated table cell to new values, in document\n// order, exactly matching the target revision.\n\nconst titleNew = \"2023-12-08 Friday\";\n\n// New values for the populated rows (row indices 0, 4, 8, 12, 16), 5 columns\n// each, in left-to-right / top-to-bottom document order.\nconst newRows = [\n  [\"76\u00f79=8, 4\", \"48\u00f78=6, 0\", \"46\u00f78=5, 6\", \"98\u00f78=12, 2\", \"56\u00f76=9, 2\"],\n  [\"62\u00f72=31, 0\", \"74\u00f73=24, 2\", \"65\u00f72=32, 1\", \"62\u00f77=8, 6\", \"42\u00f72=21, 0\"],\n  [\"53\u00f73=17, 2\", \"87\u00f73=29, 0\", \"18\u00f78=2, 2\", \"90\u00f78=11, 2\", \"13\u00f79=1, 4\"],\n  [\"32\u00f73=10, 2\", \"25\u00f76=4, 1\", \"78\u00f74=19, 2\", \"77\u00f77=11, 0\", \"72\u00f73=24, 0\"],\n  [\"85\u00f75=17, 0\", \"64\u00f74=16, 0\", \"35\u00f77=5, 0\", \"50\u00f72=25, 0\", \"78\u00f74=19, 2\"],\n];\nconst populatedRowIndices = [0, 4, 8, 12, 16];\n\nconst body = context.document.body;\n\n// --- Update the title paragraph (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(titleNew, \"Replace\");\n\n// --- Update the table cells ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let i = 0; i < populatedRowIndices.length; i++) {\n  const rowIndex = populatedRowIndices[i];\n  const rowValues = newRows[i];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.insertText(rowValues[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document is a \"two-digit division\" worksheet: a centered date title\n# paragraph followed by a 20-row x 5-column table where every 4th row\n# (1-indexed rows 1, 5, 9, 13, 17) holds a division problem/answer string\n# like \"85\u00f75=17, 0\" and the other rows are blank spacer rows. This script\n# updates the date title and every populated table cell to new values, in\n# document order, exactly matching the target revision.\n\n$d = $word.ActiveDocument\n\n# --- Update the title paragraph (first paragraph in the body) ---\n$d.Paragraphs.Item(1).Range.Text = \"2023-12-08 Friday\"\n\n# --- Update the table cells ---\n$t = $d.Tables.Item(1)\n\n# Populated 1-indexed rows, and their new left-to-right cell values.\n$rows = @(1, 5, 9, 13, 17)\n$values = @(\n    @(\"76\u00f79=8, 4\", \"48\u00f78=6, 0\", \"46\u00f78=5, 6\", \"98\u00f78=12, 2\", \"56\u00f76=9, 2\"),\n    @(\"62\u00f72=31, 0\", \"74\u00f73=24, 2\", \"65\u00f72=32, 1\", \"62\u00f77=8, 6\", \"42\u00f72=21, 0\"),\n    @(\"53\u00f73=17, 2\", \"87\u00f73=29, 0\", \"18\u00f78=2, 2\", \"90\u00f78=11, 2\", \"13\u00f79=1, 4\"),\n    @(\"32\u00f73=10, 2\", \"25\u00f76=4, 1\", \"78\u00f74=19, 2\", \"77\u00f77=11, 0\", \"72\u00f73=24, 0\"),\n    @(\"85\u00f75=17, 0\", \"64\u00f74=16, 0\", \"35\u00f77=5, 0\", \"50\u00f72=25, 0\", \"78\u00f74=19, 2\")\n)\n\nfor ($i = 0; $i -lt $rows.Count; $i++) {\n    $r = $rows[$i]\n    $rowValues = $values[$i]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
